$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain decimal number (e.g. "375.60") would be
# auto-coerced to a numeric cell by plain Range.Value assignment (losing the trailing
# zero / exact textual form). Force those through a Text number format, then clear the
# format back off again so the cell keeps its original (default) style but the value
# stays the exact literal string, just like the source inline-string cells.

$ws.Range('D2').Value = '51.187.83'
$ws.Range('E2').Value = '  -0.72%  '

$ws.Range('D3').Value = '2.941.36'
$ws.Range('E3').Value = '  -1.51%  '

$ws.Range('E4').Value = '  -0.04%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '375.60'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.65%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '102.70'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -3.08%  '

$ws.Range('E7').Value = '  -1.66%  '

$ws.Range('E8').Value = '  -0.01%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.584'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -2.49%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '36.77'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -1.88%  '

$ws.Range('E11').Value = '  -0.85%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0837'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.92%  '

$ws.Range('D13').Value = '3.401.07'
$ws.Range('E13').Value = '  -1.70%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '17.94'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -3.95%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.35'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -2.39%  '

$ws.Range('D16').Value = '2.993.87'
$ws.Range('E16').Value = '  +0.08%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.976'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +0.46%  '

$ws.Range('D18').Value = '51.115.04'
$ws.Range('E18').Value = '  -1.03%  '

$ws.Range('E19').Value = '  -7.34%  '

$ws.Range('E20').Value = '  -4.23%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '12.59'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -3.48%  '

$ws.Range('E22').Value = '  -1.09%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '263.28'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.12%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '68.26'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -1.38%  '

$ws.Range('E25').Value = '  +2.17%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '8.17'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +8.36%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.86'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +8.25%  '

$ws.Range('E28').Value = '  -2.43%  '

$ws.Range('E29').Value = '  +0.03%  '

$ws.Range('E30').Value = '  +3.57%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '25.70'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -1.43%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '9.87'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -0.34%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '34.26'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -2.44%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '50.94'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.82%  '

$ws.Range('E35').Value = '  -1.51%  '

$ws.Range('E36').Value = '  -3.68%  '

$ws.Range('E37').Value = '  -0.17%  '

$ws.Range('E38').Value = '  -3.97%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.57'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -2.08%  '

$ws.Range('B40').Value = 'Celestia'
$ws.Range('C40').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '16.46'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -5.93%  '

$ws.Range('B41').Value = 'Stellar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.115'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -1.41%  '

$ws.Range('E42').Value = '  -3.97%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '121.02'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -2.08%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '21.12'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -5.83%  '

$ws.Range('E45').Value = '  -1.66%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.273'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -3.69%  '

$ws.Range('E47').Value = '  -2.93%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '3.23'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.96%  '

$ws.Range('D49').Value = '1.999.06'
$ws.Range('E49').Value = '  -2.76%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0352'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -1.09%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '5.04'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -3.16%  '
